$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 9, shifting existing rows 9-49 down to 10-50
$ws.Rows.Item(9).Insert()

# Populate the newly inserted row 9 with the new weekly data point
$ws.Range("A9").Value = 1
$ws.Range("B9").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C9").Value = "Arica y Parinacota"
$ws.Range("D9").Value = 44459
$ws.Range("E9").Value = 15
$ws.Range("F9").Value = 100112040
$ws.Range("G9").Value = "Cilantro"
$ws.Range("H9").Value = "Sin especificar"
$ws.Range("I9").Value = "Primera"
$ws.Range("J9").Value = 250
$ws.Range("K9").Value = 900
$ws.Range("L9").Value = 1000
$ws.Range("M9").Value = 950
$ws.Range("N9").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O9").Value = "Región de Arica y Parinacota"
$ws.Range("P9").Value = 475
$ws.Range("Q9").Value = 2
$ws.Range("R9").Value = "Hortaliza"
